# previsao_retorno.xlsx - refresh of "Resumo_por_Cliente" data
# (atualizei dados bibi e add)
#
# The source data was recalculated against a newer "as of" date, which
# nudges several "INATIVO - X.Y meses sem comprar" situacao strings up by
# a tenth of a month, and refreshes the most recently active client row
# (BEMOL S/A) with its latest purchase count / dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumo_por_Cliente")

# --- situacao (column J) recalculated "months without buying" values ---
$ws.Range("J51").Value  = "INATIVO - 7.9 meses sem comprar"
$ws.Range("J66").Value  = "INATIVO - 28.2 meses sem comprar"
$ws.Range("J70").Value  = "INATIVO - 11.7 meses sem comprar"
$ws.Range("J77").Value  = "INATIVO - 8.1 meses sem comprar"
$ws.Range("J89").Value  = "INATIVO - 15.4 meses sem comprar"
$ws.Range("J91").Value  = "INATIVO - 12.1 meses sem comprar"
$ws.Range("J92").Value  = "INATIVO - 11.5 meses sem comprar"
$ws.Range("J99").Value  = "INATIVO - 33.3 meses sem comprar"
$ws.Range("J103").Value = "INATIVO - 37.6 meses sem comprar"
$ws.Range("J105").Value = "INATIVO - 14.9 meses sem comprar"

# --- row 115 (BEMOL S/A) refreshed purchase count + last/next purchase ---
$ws.Range("E115").Value = 16424
$ws.Range("H115").Value = 45846.75138888889
$ws.Range("I115").Value = 45847.75138888889
